$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.788.11'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '1.697.97'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.01'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3987'
$ws.Range('E7').Value = '  +2.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4049'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.003'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.469'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.55'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08806'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.550'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.985'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001343'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '1.673.48'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.56'
$ws.Range('E18').Value = '  -3.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07181'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.90'
$ws.Range('E20').Value = '  +4.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.343'
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.004'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.39'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').Value = '24.775.80'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.380'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.903'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.09'
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.144'
$ws.Range('E28').Value = '  +17.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '161.88'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '144.39'
$ws.Range('E30').Value = '  +5.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.353'
$ws.Range('E31').Value = '  -5.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.266'
$ws.Range('E32').Value = '  +15.44%  '
$ws.Range('D33').Value = '1.885.55'
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08634'
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.03184'
$ws.Range('E35').Value = '  +8.90%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.322'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.028'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2850'
$ws.Range('E38').Value = '  +3.56%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8320'
$ws.Range('E39').Value = '  +5.07%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09434'
$ws.Range('E40').Value = '  +3.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.71'
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.17'
$ws.Range('E42').Value = '  -0.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.476'
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('E44').Value = '  +4.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.711'
$ws.Range('E45').Value = '  +4.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7439'
$ws.Range('E46').Value = '  +3.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.220'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.379'
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.08369'
$ws.Range('E50').Value = '  +5.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '139.41'
